$d = $word.ActiveDocument

$d.Content.Find.Execute("95×81=", $true, $false, $false, $false, $false, $true, 1, $false, "47×76=", 2) | Out-Null
$d.Content.Find.Execute("37×78=", $true, $false, $false, $false, $false, $true, 1, $false, "37×79=", 2) | Out-Null
$d.Content.Find.Execute("65×66=", $true, $false, $false, $false, $false, $true, 1, $false, "14×13=", 2) | Out-Null
$d.Content.Find.Execute("46×40=", $true, $false, $false, $false, $false, $true, 1, $false, "61×25=", 2) | Out-Null
$d.Content.Find.Execute("15×98=", $true, $false, $false, $false, $false, $true, 1, $false, "11×19=", 2) | Out-Null
$d.Content.Find.Execute("19×51=", $true, $false, $false, $false, $false, $true, 1, $false, "23×51=", 2) | Out-Null
$d.Content.Find.Execute("58×53=", $true, $false, $false, $false, $false, $true, 1, $false, "18×86=", 2) | Out-Null
$d.Content.Find.Execute("64×59=", $true, $false, $false, $false, $false, $true, 1, $false, "50×42=", 2) | Out-Null
$d.Content.Find.Execute("74×71=", $true, $false, $false, $false, $false, $true, 1, $false, "13×22=", 2) | Out-Null
$d.Content.Find.Execute("82×45=", $true, $false, $false, $false, $false, $true, 1, $false, "22×52=", 2) | Out-Null
$d.Content.Find.Execute("33×66=", $true, $false, $false, $false, $false, $true, 1, $false, "67×68=", 2) | Out-Null
$d.Content.Find.Execute("61×86=", $true, $false, $false, $false, $false, $true, 1, $false, "82×75=", 2) | Out-Null
$d.Content.Find.Execute("25×13=", $true, $false, $false, $false, $false, $true, 1, $false, "28×56=", 2) | Out-Null
$d.Content.Find.Execute("80×59=", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=", 2) | Out-Null
$d.Content.Find.Execute("79×44=", $true, $false, $false, $false, $false, $true, 1, $false, "71×99=", 2) | Out-Null
$d.Content.Find.Execute("80×86=", $true, $false, $false, $false, $false, $true, 1, $false, "63×94=", 2) | Out-Null
$d.Content.Find.Execute("56×52=", $true, $false, $false, $false, $false, $true, 1, $false, "38×43=", 2) | Out-Null
$d.Content.Find.Execute("77×37=", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=", 2) | Out-Null
$d.Content.Find.Execute("14×92=", $true, $false, $false, $false, $false, $true, 1, $false, "93×45=", 2) | Out-Null
$d.Content.Find.Execute("87×16=", $true, $false, $false, $false, $false, $true, 1, $false, "49×28=", 2) | Out-Null
$d.Content.Find.Execute("28×20=", $true, $false, $false, $false, $false, $true, 1, $false, "49×40=", 2) | Out-Null
$d.Content.Find.Execute("25×55=", $true, $false, $false, $false, $false, $true, 1, $false, "13×27=", 2) | Out-Null
$d.Content.Find.Execute("19×30=", $true, $false, $false, $false, $false, $true, 1, $false, "90×62=", 2) | Out-Null
$d.Content.Find.Execute("95×12=", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=", 2) | Out-Null
$d.Content.Find.Execute("72×85=", $true, $false, $false, $false, $false, $true, 1, $false, "95×98=", 2) | Out-Null
